# Add "RIGHT" option support (like DOWN, but a different direction).
# On the "Info" sheet, add two new SPREADSHEETFORM:RIGHT blocks ("hungry"
# and "sleepy"), mirroring the existing SPREADSHEETFORM:SINGLE block
# (Pet/Noise) already on rows 5-6, and widen three extra columns (C,D,E)
# to hold the RIGHT-direction answer cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# --- Row 9/10: "hungry" RIGHT block -----------------------------------
# Column A/B formatting is copied from the existing A5/B5 label+key cells
# so the new cells share the same (bold-label / italic-key) cell styles.
$ws.Range("A5").Copy($ws.Range("A9"))
$ws.Range("B5").Copy($ws.Range("B9"))
$ws.Range("A5").Copy($ws.Range("A10"))
$ws.Range("B5").Copy($ws.Range("B10"))

$ws.Range("A9").Value = "Hungry"
$ws.Range("B9").Value = "SPREADSHEETFORM:RIGHT:hungry:state"
$ws.Range("A10").Value = "Wants"
$ws.Range("B10").Value = "SPREADSHEETFORM:RIGHT:hungry:wants"

# --- Row 13/14: "sleepy" RIGHT block ------------------------------------
$ws.Range("A5").Copy($ws.Range("A13"))
$ws.Range("B5").Copy($ws.Range("B13"))
$ws.Range("A5").Copy($ws.Range("A14"))
$ws.Range("B5").Copy($ws.Range("B14"))

$ws.Range("A13").Value = "Sleepy"
$ws.Range("B13").Value = "SPREADSHEETFORM:RIGHT:sleepy:state"
$ws.Range("A14").Value = "Wants"
$ws.Range("B14").Value = "SPREADSHEETFORM:RIGHT:sleepy:wants"

# --- Answer cells to the right of each RIGHT row (C:E), boxed with a
#     thin border like the rest of the form's input cells.
$ws.Range("C9:E10").Borders.LineStyle = 1
$ws.Range("C13:E14").Borders.LineStyle = 1

# --- Trailing blank spacer rows (16, 17), present in the sheet but with
#     no cell content - keep the form visually separated from anything
#     below it.
$ws.Rows.Item(16).RowHeight = 15
$ws.Rows.Item(17).RowHeight = 15

# --- Widen the new answer columns (C, D, E) to comfortably fit input.
# Excel's ColumnWidth is in characters; it is offset from the stored
# column width (in the xlsx <col> element) by 5/6 of a character.
$ws.Columns.Item(3).ColumnWidth = 22.62 - 0.8333333333333333
$ws.Columns.Item(4).ColumnWidth = 22.51 - 0.8333333333333333
$ws.Columns.Item(5).ColumnWidth = 23.2 - 0.8333333333333333

# --- Restore the active selection to where the form now ends.
$ws.Range("B19").Select() | Out-Null
